$d = $word.ActiveDocument

# Target the task-distribution table: row 3 ("Terence"), second column
# (the description/file cell), which currently contains a single empty
# paragraph. Fill it with the integration-test description and add a new
# paragraph below it with the corresponding test-script filename, matching
# the pattern used by the other rows in the table.
$tbl = $d.Tables.Item(1)
$cell = $tbl.Cell(3, 2)
$para = $cell.Range.Paragraphs.Item(1)

$target = $para.Range
$target.InsertAfter("Integration test for login and booking history`rLogin_booking_integration_testing.py")
$target.Font.Size = 11
